# Region.xlsx — complete the <el-option> list: duplicate a curated subset
# of the existing B:C "region" rows into a new H:I block for rows 3-18.
# (Mirrors commit "refract and complete <el-option>".)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Destination row (H/I) -> source row (B/C) to copy value+style from.
$pairs = @(
    @(3, 3),
    @(4, 6),
    @(5, 10),
    @(6, 12),
    @(7, 13),
    @(8, 14),
    @(9, 16),
    @(10, 17),
    @(11, 18),
    @(12, 21),
    @(13, 23),
    @(14, 27),
    @(15, 30),
    @(16, 37),
    @(17, 38),
    @(18, 40)
)

foreach ($pair in $pairs) {
    $dstRow = $pair[0]
    $srcRow = $pair[1]
    $ws.Range("B$srcRow`:C$srcRow").Copy($ws.Range("H$dstRow")) | Out-Null
}

# sheetView housekeeping: drop the stale scroll anchor, move the selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M18").Select() | Out-Null

# Row heights grew slightly across the whole sheet (one thicker banner row).
$ws.Rows("1:41").RowHeight = 15
$ws.Rows(28).RowHeight = 20.25
